# "upgrade left table until javakheti"
# Zestaponi municipality sheet: fix the misspelled "Zestafoni" -> "Zestaponi"
# (sheet name + the long title string) and add the 2023 data column (K)
# to the remuneration table, reusing column J's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the municipality name spelling -------------------------------
$ws.Name = "Zestaponi"
$ws.Range("A1").Value = "Average monthly remuneration of employed persons of business sector in Zestaponi Municipality"

# --- Add the new 2023 column (K), cloning column J's number formats ---
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1732.3
$ws.Range("K5").Value = 1023.5
$ws.Range("K6").Value = 1963.7
